$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2018", $false, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# 2. Fill in the empty English objectives paragraph (italic run with no text)
$objectivesText = "To update the knowledge of Biochemistry dealing with bioenergetics and metabolism. Understanding the importance of biological macromolecules and their metabolic pathways. To use this knowlege as a pre-requisite for one of the discipline of the Biochemical Engineering course"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $p.Range.Font.Italic) {
        $p.Range.Text = $objectivesText
        break
    }
}

# 3. Update "Programa resumido" Portuguese summary text
$d.Content.Find.Execute("Introdução ao metabolismo, BioenergéticaOxidações biológicas, TransporteGlicídios - metabolismoFotossínteseLipídeos - metabolismoAminoácidos - metabolismoIntegração MetabólicaCiclos vitais: oxigênio, carbono, nitrogênio e enxofre", $false, $false, $false, $false, $false, $true, 1, $false, "Introdução ao metabolismo. Conceitos básicos de metabolismo de carboidratos, gorduras e proteínas. Bioenergética, Oxidações biológicas, Transporte. Fotossíntese. Lipídeos - metabolismo. Fosforilação em nível de substrato. Aminoácidos - metabolismo. Integração Metabólica. Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre.", 2) | Out-Null

# 4. Update "Programa resumido" English summary text
$d.Content.Find.Execute("Introduction to metabolism. Bioenergetics. Biological oxidations. Glycides transport  metabolism,  photosynthesis,  lipid metabolism, metabolism of nitrogenous compounds: amino acids, integration and control of metabolic processes, vital cycles: oxygen, carbon, nitrogen and sulfur.", $false, $false, $false, $false, $false, $true, 1, $false, "Introduction to metabolism. Basic concepts of carbohydrate, fat and protein metabolism. Bioenergetics, Biological oxidations, Transport. Photosynthesis. Lipids - metabolism. Substrate-level phosphorylation. Amino acids - metabolism. Metabolic Integration. Life cycles: oxygen, carbon, nitrogen and sulfur.", 2) | Out-Null

# 5. Update "Programa" Portuguese full syllabus text
$d.Content.Find.Execute("Introdução ao metabolismo. Proteínas, polissacarídeos, lipídios: vias catabólicas e anabólicas. Bioenergética. Variação de energia livre: relação com a constante de equilíbrio e com o potencial redox. Processos exergônicos. Papel do fosfato: potencial de transferência de grupo fosfato. Importância energética do ATP.Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa.Transporte. Composição das membranas biológicas: constituição química, caráter ?barreira permeabilidade?. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: glicólise - localização das enzimas operantes, reações, produção de NADH.H+, fosforilação ao nível de substrato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. ?Shunt? das hexoses-fosfato (ciclo das pentoses): localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH.H  (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP . Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin.Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da b-oxidação. a-oxidação e w-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos.Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos.Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos.Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato", $false, $false, $false, $false, $false, $true, 1, $false, "Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa. fosforilação em nível de substratoTransporte. Composição das membranas biológicas: constituição química, caráter barreira permeabilidade. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: básico e conceito de reações enzimáticas. reações, produção de NADH.H+, descarboxilação de piruvato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. Via das pentoses-fosfato: localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH. (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica.  Conceitos básicos de fermentação em estado sólido, fermentação semissólida e fermentação submersa. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP. Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin. Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da β-oxidação. α-oxidação e ω-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos. Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos. Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos. Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato.", 2) | Out-Null

# 6. Update "Programa" English full syllabus text
$d.Content.Find.Execute("Introduction to metabolism. Proteins, polysaccharides, lipids.: catabolic and anabolic pathways. Bioenergetics. The direction of processes: free energy: reaction with the balance constant and with redox potential. Exergonic processes. Role of phosphate: potential of transference of phosphate group. Role of ATP as a free energy currency. Biological oxidations. Electron and proton transporters coenzymes: nucleotides, flavin nucleotides, coenzyme Q. Pyridine and flavin nucleotides- dehydrogenase dependent. Oxidases.Structure of mitochondrial membrane. Respiratory chain: function. Oxidative phosphorylation. Transport. Composition of biological membranes: chemical constitution, characterization, barrier, permeability. Carriers and canal ionophores. Transportation: mediated and non-mediated. Glycides  metabolism. Anaerobic and aerobic degradation of glycides: glycolysis  localization of enzymes, reactions, NaDH.H+ production, the first substrt level phosphorylation, energetic balance; Citric acid cycle  localization of operating enzymes, reactions, production of reduced coenzymes, energetic balance. the pentose phosphate pathway: localization of the enzymes, reactions (oxidative and non-oxidative phases), NaDH.H+ production (physiological implication). Fermentation: definition, fermentation and respiration, raw-materials used in the starch and sugar fermentations, metabolic fates of pyruvate: ethanol and lactic metabolism, acetic and citric. Photosynthesis. Chloroplasts structure. The light reactions. Role of the chlorophyll in the photosynthesis. Cyclic and non-cyclic phosphorylation. NADP reduction. Water photolysis. Synthesis of the acceptor of CO2, Ru-1, 5-diP. Calvin cycle. Lipids  metabolism. Beta-oxidation pathways, oxidation of fatty acids with odd-numbered carbon chains, energetic balance of beta-oxidation, alpha-oxidation and w-oxidation.  Glycerol metabolism. Formation of ketone bodies. Biosynthesis of fatty acids. Amino acids  catabolic pathways. Digestion of proteins, aspects of amino acid synthesis and degradation.Transamination, urea cycle. Metabolic integration. Common metabolites to the metabolism of glycides, lipids and amino acids. Vital cycles: oxygen, carbon, nitrogen and sulfur. Carbon and oxygen cycles. Nitrogen cycle: biological fixation, nitrification, use of nitrate, incorporation of ammonia in organic compounds. Sulfur cycle: sulfate assimilation.", $false, $false, $false, $false, $false, $true, 1, $false, "Biological oxidation. Proton and electron transport coenzymes: nucleotides, flavin nucleotides, coenzyme-Q. Pyridine and flavin nucleotide-dependent dehydrogenase. Oxidases. Structure of the mitochondrial membrane. Respiratory chain: function. Oxidative phosphorylation. substrate-level phosphorylationTransport. Composition of biological membranes: chemical constitution, permeability barrier character. Carriers and channels, ionophores. Transport: mediated and non-mediated. Carbohydrates - metabolism. Anaerobic and aerobic degradation of carbohydrates: basics and concept of enzymatic reactions. reactions, production of NADH.H+, decarboxylation of pyruvate, energy balance; Krebs cycle - location of operating enzymes, reactions, production of reduced coenzymes, energy balance. Hexose-phosphate pathway (pentose cycle): location of operating enzymes, reactions (oxidative and non-oxidative phases), production of NADPH. (physiological implication). Fermentations: definition, difference between fermentation and respiration, raw materials used in starchy and sugars fermentation, fermentation agent, anaerobic fermentations - alcoholic and lactic, aerobic fermentations - acetic and citric. Basic concepts of solid state fermentation, sem-isolid fermentation and submerged fermentation. Photosynthesis. Structure of chloroplasts. Light dependente reactions. Role of chlorophyll in photosynthesis. Cyclic and non-cyclic photophosphorylation. NADP reduction. Photolysis of water. Synthesis of the CO2 acceptor, Calvin cycle. Lipids - metabolism. β-oxidation of chain fatty acids: with an even number of C atoms, with an odd number of C atoms, branched, energy balance of β-oxidation. α-oxidation and ω-oxidation. Glycerol metabolism. Formation of ketone bodies. Biosynthesis of fatty acids. Amino acids - catabolic pathways. Digestion of proteins, involved enzymes and zymogens, absorption. Transamination, oxidative deamination, ketogenic and glucogenic amino acids, decarboxylation, urea cycle. Nitrogen elimination, ureotelic vertebrates. Metabolic integration. Interrelation of the intermediate metabolism of carbohydrates, lipids, amino acids and nucleic acids. Life cycles: oxygen, carbon, nitrogen and sulfur. Oxygen and carbon cycles. Nitrogen cycle: Biological fixation, nitrification, use of nitrate, incorporation of ammonia into organic compounds. Sulfur cycle: Assimilation of sulfate", 2) | Out-Null

# 7. Update the final grade formula text (Criterio)
$d.Content.Find.Execute("A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2", $false, $false, $false, $false, $false, $true, 1, $false, "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2.", 2) | Out-Null

